# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets to
# reflect the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4224
$ws1.Range("F3").Value = 2406
$ws1.Range("F5").Value = 19
$ws1.Range("F7").Value = 45
$ws1.Range("F10").Value = 121
$ws1.Range("F12").Value = 1578
$ws1.Range("F13").Value = 290
$ws1.Range("F14").Value = 3231
$ws1.Range("F15").Value = 219

# --- 演出 sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 42

# --- 全部类型 sheet -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4224
$ws4.Range("F3").Value = 2406
$ws4.Range("F5").Value = 19
$ws4.Range("F8").Value = 45
$ws4.Range("F9").Value = 42
$ws4.Range("F12").Value = 121
$ws4.Range("F16").Value = 1578
$ws4.Range("F17").Value = 290
$ws4.Range("F18").Value = 3231
$ws4.Range("F19").Value = 219
